# Insert three new weekly price rows at the top of the Chirimoya data block
# (rows 66-68), shifting the existing rows 66-137 down to 69-140.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 66 (pushes existing data down by 3 rows)
$ws.Range("A66:A68").EntireRow.Insert()

# Common (constant) field values shared by every Chirimoya / Comercializadora
# del Agro de Limarí record in this sheet.
$mercadoId = 2
$mercado   = "Comercializadora del Agro de Limarí"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$origen    = "Provincia de Limarí"
$fecha     = 44902

# Row 66: Especial
$ws.Range("A66").Value = $mercadoId
$ws.Range("B66").Value = $mercado
$ws.Range("C66").Value = $region
$ws.Range("D66").Value = $fecha
$ws.Range("E66").Value = $codreg
$ws.Range("F66").Value = $tipo
$ws.Range("G66").Value = $productoId
$ws.Range("H66").Value = $producto
$ws.Range("I66").Value = $categoriaId
$ws.Range("J66").Value = $categoria
$ws.Range("K66").Value = $variedad
$ws.Range("L66").Value = "Especial"
$ws.Range("M66").Value = 200
$ws.Range("N66").Value = 15000
$ws.Range("O66").Value = 16000
$ws.Range("P66").Value = 15500
$ws.Range("Q66").Value = '$/bandeja 10 kilos'
$ws.Range("R66").Value = $origen
$ws.Range("S66").Value = 1550
$ws.Range("T66").Value = 10

# Row 67: Primera
$ws.Range("A67").Value = $mercadoId
$ws.Range("B67").Value = $mercado
$ws.Range("C67").Value = $region
$ws.Range("D67").Value = $fecha
$ws.Range("E67").Value = $codreg
$ws.Range("F67").Value = $tipo
$ws.Range("G67").Value = $productoId
$ws.Range("H67").Value = $producto
$ws.Range("I67").Value = $categoriaId
$ws.Range("J67").Value = $categoria
$ws.Range("K67").Value = $variedad
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 240
$ws.Range("N67").Value = 13000
$ws.Range("O67").Value = 14000
$ws.Range("P67").Value = 13500
$ws.Range("Q67").Value = '$/bandeja 10 kilos'
$ws.Range("R67").Value = $origen
$ws.Range("S67").Value = 1350
$ws.Range("T67").Value = 10

# Row 68: Segunda
$ws.Range("A68").Value = $mercadoId
$ws.Range("B68").Value = $mercado
$ws.Range("C68").Value = $region
$ws.Range("D68").Value = $fecha
$ws.Range("E68").Value = $codreg
$ws.Range("F68").Value = $tipo
$ws.Range("G68").Value = $productoId
$ws.Range("H68").Value = $producto
$ws.Range("I68").Value = $categoriaId
$ws.Range("J68").Value = $categoria
$ws.Range("K68").Value = $variedad
$ws.Range("L68").Value = "Segunda"
$ws.Range("M68").Value = 300
$ws.Range("N68").Value = 9000
$ws.Range("O68").Value = 10000
$ws.Range("P68").Value = 9500
$ws.Range("Q68").Value = '$/bandeja 10 kilos'
$ws.Range("R68").Value = $origen
$ws.Range("S68").Value = 950
$ws.Range("T68").Value = 10
